# Auto-generated edit script: applies cell-value changes per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value = 1139.8
$ws.Cells.Item(5, 9).Value = 170
$ws.Cells.Item(5, 10).Value = 2594.5
$ws.Cells.Item(5, 11).Value = 170
$ws.Cells.Item(5, 12).Value = 2594.5
$ws.Cells.Item(5, 13).Value = -55
$ws.Cells.Item(5, 14).Value = -2824.5

$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 9).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 11).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 13).Value = ""
$ws.Cells.Item(32, 14).Value = ""

$ws.Cells.Item(43, 8).Value = 8018.3335
$ws.Cells.Item(43, 9).Value = 0
$ws.Cells.Item(43, 10).Value = 8018.3335
$ws.Cells.Item(43, 11).Value = 0
$ws.Cells.Item(43, 12).Value = 8018.3335
$ws.Cells.Item(43, 13).Value = ""
$ws.Cells.Item(43, 14).Value = -8156.3335

$ws.Cells.Item(76, 8).Value = 6004
$ws.Cells.Item(76, 10).Value = 6004
$ws.Cells.Item(76, 12).Value = 6004
$ws.Cells.Item(76, 14).Value = -6634

$ws.Cells.Item(79, 8).Value = 6004
$ws.Cells.Item(79, 10).Value = 6004
$ws.Cells.Item(79, 12).Value = 6004
$ws.Cells.Item(79, 14).Value = -8188

$ws.Cells.Item(86, 8).Value = 66698984
$ws.Cells.Item(86, 9).Value = 1639.7778
$ws.Cells.Item(86, 11).Value = 1639.7778
$ws.Cells.Item(86, 13).Value = -516.7778000000001

$ws.Cells.Item(87, 8).Value = 67992.57000000001
$ws.Cells.Item(87, 10).Value = 67992.57000000001
$ws.Cells.Item(87, 12).Value = 67992.57000000001
$ws.Cells.Item(87, 14).Value = -70488.57000000001

$ws.Cells.Item(89, 8).Value = 66698984
$ws.Cells.Item(89, 9).Value = 1639.7778
$ws.Cells.Item(89, 11).Value = 8198.889000000001
$ws.Cells.Item(89, 13).Value = -2582.889000000001

$ws.Cells.Item(90, 8).Value = 67992.57000000001
$ws.Cells.Item(90, 10).Value = 67992.57000000001
$ws.Cells.Item(90, 12).Value = 203977.71
$ws.Cells.Item(90, 14).Value = -216457.71

$ws.Cells.Item(101, 8).Value = 15877009
$ws.Cells.Item(101, 9).Value = 17861386
$ws.Cells.Item(101, 10).Value = 1998
$ws.Cells.Item(101, 11).Value = 53584158
$ws.Cells.Item(101, 12).Value = 5994
$ws.Cells.Item(101, 13).Value = -53582536
$ws.Cells.Item(101, 14).Value = -9238

$ws.Cells.Item(109, 8).Value = 105278350
$ws.Cells.Item(109, 10).Value = 105278350
$ws.Cells.Item(109, 12).Value = 105278350
$ws.Cells.Item(109, 14).Value = -105281124

$ws.Cells.Item(113, 8).Value = 8403.619000000001
$ws.Cells.Item(113, 9).Value = 7915.5835
$ws.Cells.Item(113, 11).Value = 7915.5835
$ws.Cells.Item(113, 13).Value = -4661.5835

$ws.Cells.Item(131, 8).Value = 5131.125
$ws.Cells.Item(131, 9).Value = 3409.8
$ws.Cells.Item(131, 11).Value = 10229.4
$ws.Cells.Item(131, 13).Value = -5189.400000000001

$ws.Cells.Item(132, 8).Value = 4032.7532
$ws.Cells.Item(132, 9).Value = 4367.9
$ws.Cells.Item(132, 11).Value = 13103.7
$ws.Cells.Item(132, 13).Value = -10573.7

$ws.Cells.Item(138, 8).Value = 270948.62
$ws.Cells.Item(138, 10).Value = 4046.8125
$ws.Cells.Item(138, 12).Value = 12140.4375
$ws.Cells.Item(138, 14).Value = -22420.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(44, 8).Value = 70030.5
$ws.Cells.Item(44, 10).Value = 70030.5
$ws.Cells.Item(44, 12).Value = 70030.5
$ws.Cells.Item(44, 14).Value = -71006.5

$ws.Cells.Item(45, 8).Value = 295571.28
$ws.Cells.Item(45, 9).Value = 508999.75
$ws.Cells.Item(45, 10).Value = 11000
$ws.Cells.Item(45, 11).Value = 508999.75
$ws.Cells.Item(45, 12).Value = 11000
$ws.Cells.Item(45, 13).Value = -508622.75
$ws.Cells.Item(45, 14).Value = -11754

$ws.Cells.Item(61, 8).Value = 3633.246
$ws.Cells.Item(61, 9).Value = 3357.32
$ws.Cells.Item(61, 10).Value = 4553
$ws.Cells.Item(61, 11).Value = 3357.32
$ws.Cells.Item(61, 12).Value = 4553
$ws.Cells.Item(61, 13).Value = -3145.32
$ws.Cells.Item(61, 14).Value = -4977

$ws.Cells.Item(74, 8).Value = 4054.3333
$ws.Cells.Item(74, 9).Value = 2332.9565
$ws.Cells.Item(74, 11).Value = 2332.9565
$ws.Cells.Item(74, 13).Value = -1458.9565

$ws.Cells.Item(77, 8).Value = 4054.3333
$ws.Cells.Item(77, 9).Value = 2332.9565
$ws.Cells.Item(77, 11).Value = 11664.7825
$ws.Cells.Item(77, 13).Value = -7296.782499999999

$ws.Cells.Item(122, 8).Value = 1423740.6
$ws.Cells.Item(122, 10).Value = 3389487
$ws.Cells.Item(122, 12).Value = 10168461
$ws.Cells.Item(122, 14).Value = -10173361

$ws.Cells.Item(132, 8).Value = 12156
$ws.Cells.Item(132, 9).Value = 15225.444
$ws.Cells.Item(132, 11).Value = 45676.33199999999
$ws.Cells.Item(132, 13).Value = -43146.33199999999

$ws.Cells.Item(136, 8).Value = 3633.246
$ws.Cells.Item(136, 9).Value = 3357.32
$ws.Cells.Item(136, 10).Value = 4553
$ws.Cells.Item(136, 11).Value = 10071.96
$ws.Cells.Item(136, 12).Value = 13659
$ws.Cells.Item(136, 13).Value = -7521.960000000001
$ws.Cells.Item(136, 14).Value = -18759

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 3202.5625
$ws.Cells.Item(20, 9).Value = 2907.3333
$ws.Cells.Item(20, 10).Value = 3379.7
$ws.Cells.Item(20, 11).Value = 2907.3333
$ws.Cells.Item(20, 12).Value = 3379.7
$ws.Cells.Item(20, 13).Value = -2660.3333
$ws.Cells.Item(20, 14).Value = -3873.7

$ws.Cells.Item(134, 8).Value = 11312.186
$ws.Cells.Item(134, 9).Value = 14003
$ws.Cells.Item(134, 11).Value = 42009
$ws.Cells.Item(134, 13).Value = -39474

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4530.3335
$ws.Cells.Item(31, 9).Value = 4046.32
$ws.Cells.Item(31, 10).Value = 6042.875
$ws.Cells.Item(31, 11).Value = 4046.32
$ws.Cells.Item(31, 12).Value = 6042.875
$ws.Cells.Item(31, 13).Value = -3751.32
$ws.Cells.Item(31, 14).Value = -6632.875

$ws.Cells.Item(34, 8).Value = 4530.3335
$ws.Cells.Item(34, 9).Value = 4046.32
$ws.Cells.Item(34, 10).Value = 6042.875
$ws.Cells.Item(34, 11).Value = 4046.32
$ws.Cells.Item(34, 12).Value = 6042.875
$ws.Cells.Item(34, 13).Value = -3844.32
$ws.Cells.Item(34, 14).Value = -6446.875

$ws.Cells.Item(94, 8).Value = 2564.6
$ws.Cells.Item(94, 9).Value = 1875.3334
$ws.Cells.Item(94, 11).Value = 1875.3334
$ws.Cells.Item(94, 13).Value = -1424.3334

$ws.Cells.Item(132, 8).Value = 11819.7
$ws.Cells.Item(132, 9).Value = 1289.5
$ws.Cells.Item(132, 11).Value = 3868.5
$ws.Cells.Item(132, 13).Value = -1338.5

$ws.Cells.Item(134, 8).Value = 2678.4722
$ws.Cells.Item(134, 9).Value = 2204.2307
$ws.Cells.Item(134, 11).Value = 6612.6921
$ws.Cells.Item(134, 13).Value = -4077.6921

$ws.Cells.Item(135, 8).Value = 69999
$ws.Cells.Item(135, 10).Value = 69999
$ws.Cells.Item(135, 12).Value = 69999
$ws.Cells.Item(135, 14).Value = -80139

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 89619510
$ws.Cells.Item(4, 9).Value = 86984936
$ws.Cells.Item(4, 11).Value = 260954808
$ws.Cells.Item(4, 13).Value = -260954696

$ws.Cells.Item(12, 8).Value = 199.08333
$ws.Cells.Item(12, 10).Value = 138.71428
$ws.Cells.Item(12, 12).Value = 416.14284
$ws.Cells.Item(12, 14).Value = -762.14284

$ws.Cells.Item(32, 8).Value = 694.75
$ws.Cells.Item(32, 10).Value = 694.75
$ws.Cells.Item(32, 12).Value = 2084.25
$ws.Cells.Item(32, 14).Value = -2650.25

$ws.Cells.Item(46, 8).Value = 2515.75
$ws.Cells.Item(46, 10).Value = 3058.3333
$ws.Cells.Item(46, 12).Value = 9174.999899999999
$ws.Cells.Item(46, 14).Value = -9356.999899999999

$ws.Cells.Item(50, 8).Value = 671.1
$ws.Cells.Item(50, 9).Value = 301.33334
$ws.Cells.Item(50, 11).Value = 904.0000200000001
$ws.Cells.Item(50, 13).Value = -423.0000200000001

$ws.Cells.Item(53, 8).Value = 671.1
$ws.Cells.Item(53, 9).Value = 301.33334
$ws.Cells.Item(53, 11).Value = 904.0000200000001
$ws.Cells.Item(53, 13).Value = -423.0000200000001

$ws.Cells.Item(64, 8).Value = 1967.45
$ws.Cells.Item(64, 10).Value = 1788.7142
$ws.Cells.Item(64, 12).Value = 5366.142599999999
$ws.Cells.Item(64, 14).Value = -5906.142599999999

$ws.Cells.Item(67, 8).Value = 1967.45
$ws.Cells.Item(67, 10).Value = 1788.7142
$ws.Cells.Item(67, 12).Value = 5366.142599999999
$ws.Cells.Item(67, 14).Value = -7238.142599999999

$ws.Cells.Item(68, 8).Value = 22734096
$ws.Cells.Item(68, 9).Value = 1587.4
$ws.Cells.Item(68, 11).Value = 4762.200000000001
$ws.Cells.Item(68, 13).Value = -3951.200000000001

$ws.Cells.Item(70, 8).Value = 1623.3334
$ws.Cells.Item(70, 9).Value = 515.7143
$ws.Cells.Item(70, 11).Value = 1547.1429
$ws.Cells.Item(70, 13).Value = -1232.1429

$ws.Cells.Item(71, 8).Value = 22734096
$ws.Cells.Item(71, 9).Value = 1587.4
$ws.Cells.Item(71, 11).Value = 14286.6
$ws.Cells.Item(71, 13).Value = -10230.6

$ws.Cells.Item(73, 8).Value = 1623.3334
$ws.Cells.Item(73, 9).Value = 515.7143
$ws.Cells.Item(73, 11).Value = 1547.1429
$ws.Cells.Item(73, 13).Value = -455.1428999999998

$ws.Cells.Item(75, 8).Value = 464.66666
$ws.Cells.Item(75, 10).Value = 500
$ws.Cells.Item(75, 12).Value = 1500
$ws.Cells.Item(75, 14).Value = -3496

$ws.Cells.Item(78, 8).Value = 464.66666
$ws.Cells.Item(78, 10).Value = 500
$ws.Cells.Item(78, 12).Value = 4500
$ws.Cells.Item(78, 14).Value = -14484

$ws.Cells.Item(98, 8).Value = 962.7222
$ws.Cells.Item(98, 9).Value = 1138.5385
$ws.Cells.Item(98, 11).Value = 3415.6155
$ws.Cells.Item(98, 13).Value = -1917.6155

$ws.Cells.Item(113, 8).Value = 1186.359
$ws.Cells.Item(113, 10).Value = 1400.826
$ws.Cells.Item(113, 12).Value = 4202.478
$ws.Cells.Item(113, 14).Value = -8542.477999999999

$ws.Cells.Item(131, 8).Value = 34486548
$ws.Cells.Item(131, 10).Value = 2665.889
$ws.Cells.Item(131, 12).Value = 7997.667
$ws.Cells.Item(131, 14).Value = -18077.667

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(44, 8).Value = 29995
$ws.Cells.Item(44, 10).Value = 29995
$ws.Cells.Item(44, 12).Value = 29995
$ws.Cells.Item(44, 14).Value = -31187

$ws.Cells.Item(70, 8).Value = 5856
$ws.Cells.Item(70, 9).Value = 4807.125
$ws.Cells.Item(70, 11).Value = 4807.125
$ws.Cells.Item(70, 13).Value = -4537.125

$ws.Cells.Item(73, 8).Value = 5856
$ws.Cells.Item(73, 9).Value = 4807.125
$ws.Cells.Item(73, 11).Value = 4807.125
$ws.Cells.Item(73, 13).Value = -3871.125

$ws.Cells.Item(80, 8).Value = 14200.4
$ws.Cells.Item(80, 9).Value = 16125.5
$ws.Cells.Item(80, 10).Value = 6500
$ws.Cells.Item(80, 11).Value = 16125.5
$ws.Cells.Item(80, 12).Value = 6500
$ws.Cells.Item(80, 13).Value = -15127.5
$ws.Cells.Item(80, 14).Value = -8496

$ws.Cells.Item(83, 8).Value = 14200.4
$ws.Cells.Item(83, 9).Value = 16125.5
$ws.Cells.Item(83, 10).Value = 6500
$ws.Cells.Item(83, 11).Value = 80627.5
$ws.Cells.Item(83, 12).Value = 32500
$ws.Cells.Item(83, 13).Value = -75635.5
$ws.Cells.Item(83, 14).Value = -42484

$ws.Cells.Item(102, 8).Value = 4963.3877
$ws.Cells.Item(102, 9).Value = 5346.6587
$ws.Cells.Item(102, 11).Value = 5346.6587
$ws.Cells.Item(102, 13).Value = -3724.6587

$ws.Cells.Item(107, 8).Value = 587.9048
$ws.Cells.Item(107, 10).Value = 400.83334
$ws.Cells.Item(107, 12).Value = 400.83334
$ws.Cells.Item(107, 14).Value = -4240.83334

$ws.Cells.Item(122, 8).Value = 8322.629999999999
$ws.Cells.Item(122, 9).Value = 5669.864
$ws.Cells.Item(122, 10).Value = 19994.8
$ws.Cells.Item(122, 11).Value = 17009.592
$ws.Cells.Item(122, 12).Value = 59984.39999999999
$ws.Cells.Item(122, 13).Value = -14559.592
$ws.Cells.Item(122, 14).Value = -64884.39999999999

$ws.Cells.Item(126, 8).Value = 16655.559
$ws.Cells.Item(126, 9).Value = 17651.334
$ws.Cells.Item(126, 10).Value = 14265.7
$ws.Cells.Item(126, 11).Value = 52954.00199999999
$ws.Cells.Item(126, 12).Value = 42797.10000000001
$ws.Cells.Item(126, 13).Value = -50484.00199999999
$ws.Cells.Item(126, 14).Value = -47737.10000000001

$ws.Cells.Item(132, 8).Value = 4060.2083
$ws.Cells.Item(132, 9).Value = 3876.4688
$ws.Cells.Item(132, 10).Value = 4427.6875
$ws.Cells.Item(132, 11).Value = 11629.4064
$ws.Cells.Item(132, 12).Value = 13283.0625
$ws.Cells.Item(132, 13).Value = -9099.4064
$ws.Cells.Item(132, 14).Value = -18343.0625

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 12397.464
$ws.Cells.Item(22, 9).Value = 16289.353
$ws.Cells.Item(22, 10).Value = 6382.727
$ws.Cells.Item(22, 11).Value = 16289.353
$ws.Cells.Item(22, 12).Value = 6382.727
$ws.Cells.Item(22, 13).Value = -15994.353
$ws.Cells.Item(22, 14).Value = -6972.727

$ws.Cells.Item(27, 8).Value = 12397.464
$ws.Cells.Item(27, 9).Value = 16289.353
$ws.Cells.Item(27, 10).Value = 6382.727
$ws.Cells.Item(27, 11).Value = 16289.353
$ws.Cells.Item(27, 12).Value = 6382.727
$ws.Cells.Item(27, 13).Value = -16182.353
$ws.Cells.Item(27, 14).Value = -6596.727

$ws.Cells.Item(61, 8).Value = 3504.9714
$ws.Cells.Item(61, 9).Value = 2199.1428
$ws.Cells.Item(61, 10).Value = 8728.286
$ws.Cells.Item(61, 11).Value = 2199.1428
$ws.Cells.Item(61, 12).Value = 8728.286
$ws.Cells.Item(61, 13).Value = -1997.1428
$ws.Cells.Item(61, 14).Value = -9132.286

$ws.Cells.Item(113, 8).Value = 3504.9714
$ws.Cells.Item(113, 9).Value = 2199.1428
$ws.Cells.Item(113, 10).Value = 8728.286
$ws.Cells.Item(113, 11).Value = 2199.1428
$ws.Cells.Item(113, 12).Value = 8728.286
$ws.Cells.Item(113, 13).Value = -29.14280000000008
$ws.Cells.Item(113, 14).Value = -13068.286

$ws.Cells.Item(132, 8).Value = 418158.78
$ws.Cells.Item(132, 9).Value = 712401.0600000001
$ws.Cells.Item(132, 11).Value = 2137203.18
$ws.Cells.Item(132, 13).Value = -2134673.18

$ws.Cells.Item(136, 8).Value = 5624.407
$ws.Cells.Item(136, 9).Value = 3010.6
$ws.Cells.Item(136, 10).Value = 8891.666999999999
$ws.Cells.Item(136, 11).Value = 9031.799999999999
$ws.Cells.Item(136, 12).Value = 26675.001
$ws.Cells.Item(136, 13).Value = -6481.799999999999
$ws.Cells.Item(136, 14).Value = -31775.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(4, 8).Value = 1332
$ws.Cells.Item(4, 9).Value = 1164.9231
$ws.Cells.Item(4, 10).Value = 1875
$ws.Cells.Item(4, 11).Value = 1164.9231
$ws.Cells.Item(4, 12).Value = 1875
$ws.Cells.Item(4, 13).Value = -1051.9231
$ws.Cells.Item(4, 14).Value = -2101

$ws.Cells.Item(32, 8).Value = 14587.333
$ws.Cells.Item(32, 9).Value = 14587.333
$ws.Cells.Item(32, 11).Value = 14587.333
$ws.Cells.Item(32, 13).Value = -14270.333

$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 13).Value = ""

$ws.Cells.Item(81, 8).Value = 13584.846
$ws.Cells.Item(81, 9).Value = 13966.917
$ws.Cells.Item(81, 10).Value = 9000
$ws.Cells.Item(81, 11).Value = 27933.834
$ws.Cells.Item(81, 12).Value = 18000
$ws.Cells.Item(81, 13).Value = -26872.834
$ws.Cells.Item(81, 14).Value = -20122

$ws.Cells.Item(84, 8).Value = 13584.846
$ws.Cells.Item(84, 9).Value = 13966.917
$ws.Cells.Item(84, 10).Value = 9000
$ws.Cells.Item(84, 11).Value = 139669.17
$ws.Cells.Item(84, 12).Value = 90000
$ws.Cells.Item(84, 13).Value = -134365.17
$ws.Cells.Item(84, 14).Value = -100608

$ws.Cells.Item(113, 8).Value = 998.125
$ws.Cells.Item(113, 9).Value = 539.60974
$ws.Cells.Item(113, 11).Value = 1618.82922
$ws.Cells.Item(113, 13).Value = 551.1707799999999

$ws.Cells.Item(114, 8).Value = 25000
$ws.Cells.Item(114, 10).Value = 25000
$ws.Cells.Item(114, 12).Value = 25000
$ws.Cells.Item(114, 14).Value = -33678

$ws.Cells.Item(126, 8).Value = 21746.38
$ws.Cells.Item(126, 9).Value = 30892.072
$ws.Cells.Item(126, 10).Value = 3455
$ws.Cells.Item(126, 11).Value = 92676.216
$ws.Cells.Item(126, 12).Value = 10365
$ws.Cells.Item(126, 13).Value = -90206.216
$ws.Cells.Item(126, 14).Value = -15305

$ws.Cells.Item(132, 8).Value = 9917.305
$ws.Cells.Item(132, 9).Value = 11016.357
$ws.Cells.Item(132, 11).Value = 33049.071
$ws.Cells.Item(132, 13).Value = -30519.071

$ws.Cells.Item(136, 8).Value = 368231.22
$ws.Cells.Item(136, 9).Value = 428658.66
$ws.Cells.Item(136, 11).Value = 1285975.98
$ws.Cells.Item(136, 13).Value = -1283425.98

